# The sheet originally had its data starting at row 4 (rows 1-3 were blank).
# This edit removes those 3 leading blank rows, shifting every row of data
# up by three (row N -> row N-3), and updates the selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three blank leading rows - everything below shifts up by 3.
[void]$ws.Rows("1:3").Delete()

# Update the selection to A20:B20, as recorded after the row deletion.
[void]$ws.Range("A20:B20").Select()
